# fix: revert admin dev default; seed customers only when table empty; autosave on customer select when hours/day present
#
# This script applies the following corrections to the weekly timesheet
# export for Chris Zavesky (2026-01-12 week):
#   - Client names on rows 2-6 are corrected (dev/seed placeholder clients
#     swapped for the real customer names).
#   - Hours/Rate/Total are reset to reflect that the admin-dev default
#     hours seed was reverted (no hours logged yet for these rows), so
#     Rate/Total go to 0 and the first row's Hours corrects from 12 to 8.
#   - The subtotal/summary rows recompute accordingly (Reg hours 44 -> 40,
#     and the dependent Total figures drop to 0).
#   - The Employee ID is refreshed to the new seeded id format.
#   - The per-row "Notes" seed text (sample-data placeholders) is cleared
#     now that the row values are no longer seeded placeholders.

$wb = $excel.ActiveWorkbook

$wsTime   = $wb.Worksheets.Item("Weekly Timesheet")
$wsSchema = $wb.Worksheets.Item("Jason Schema")

# ---- Weekly Timesheet sheet ------------------------------------------------

# Row 2: 2026-01-12
$wsTime.Range("B2").Value = "Tormey"
$wsTime.Range("C2").Value = 8
$wsTime.Range("E2").Value = 0
$wsTime.Range("F2").Value = 0

# Row 3: 2026-01-13
$wsTime.Range("B3").Value = "Richer"
$wsTime.Range("E3").Value = 0
$wsTime.Range("F3").Value = 0

# Row 4: 2026-01-14
$wsTime.Range("B4").Value = "Durfee"
$wsTime.Range("E4").Value = 0
$wsTime.Range("F4").Value = 0

# Row 5: 2026-01-15
$wsTime.Range("B5").Value = "Tercek"
$wsTime.Range("E5").Value = 0
$wsTime.Range("F5").Value = 0

# Row 6: 2026-01-16
$wsTime.Range("B6").Value = "Patton"
$wsTime.Range("E6").Value = 0
$wsTime.Range("F6").Value = 0

# SUBTOTAL row
$wsTime.Range("C8").Value = 40
$wsTime.Range("D8").Value = "Reg: 40 / OT: 0"
$wsTime.Range("F8").Value = 0

# ADMIN SUBTOTAL / GRAND TOTAL rows
$wsTime.Range("F12").Value = 0
$wsTime.Range("F13").Value = 0

# ---- Jason Schema sheet ----------------------------------------------------

# Employee ID refreshed for all data rows
$wsSchema.Range("B2").Value = "emp_5chpvt65"
$wsSchema.Range("B3").Value = "emp_5chpvt65"
$wsSchema.Range("B4").Value = "emp_5chpvt65"
$wsSchema.Range("B5").Value = "emp_5chpvt65"
$wsSchema.Range("B6").Value = "emp_5chpvt65"

# Row 2
$wsSchema.Range("D2").Value = "Tormey"
$wsSchema.Range("E2").Value = 8
$wsSchema.Range("F2").Value = 0
$wsSchema.Range("G2").Value = 0
$wsSchema.Range("I2").Value = ""

# Row 3
$wsSchema.Range("D3").Value = "Richer"
$wsSchema.Range("F3").Value = 0
$wsSchema.Range("G3").Value = 0
$wsSchema.Range("I3").Value = ""

# Row 4
$wsSchema.Range("D4").Value = "Durfee"
$wsSchema.Range("F4").Value = 0
$wsSchema.Range("G4").Value = 0
$wsSchema.Range("I4").Value = ""

# Row 5
$wsSchema.Range("D5").Value = "Tercek"
$wsSchema.Range("F5").Value = 0
$wsSchema.Range("G5").Value = 0
$wsSchema.Range("I5").Value = ""

# Row 6
$wsSchema.Range("D6").Value = "Patton"
$wsSchema.Range("F6").Value = 0
$wsSchema.Range("G6").Value = 0
$wsSchema.Range("I6").Value = ""

$wb.Save()
